$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 180-181, pushing the existing 180-189 block down to 182-191.
$ws.Range("A180:T181").EntireRow.Insert()

# New row 180: Blue Giant / Primera
$ws.Range("A180").Value = 10
$ws.Range("B180").Value = "Vega Modelo de Temuco"
$ws.Range("C180").Value = "La Araucanía"
$ws.Range("D180").Value = 44610
$ws.Range("E180").Value = 9
$ws.Range("F180").Value = "Fruta"
$ws.Range("G180").Value = 100103
$ws.Range("H180").Value = "Frutos de hueso (carozo)"
$ws.Range("I180").Value = 100103002
$ws.Range("J180").Value = "Ciruela"
$ws.Range("K180").Value = "Blue Giant"
$ws.Range("L180").Value = "Primera"
$ws.Range("M180").Value = 100
$ws.Range("N180").Value = 13000
$ws.Range("O180").Value = 13000
$ws.Range("P180").Value = 13000
$ws.Range("Q180").Value = "$/bandeja 18 kilos granel"
$ws.Range("R180").Value = "Región de O'Higgins"
$ws.Range("S180").Value = 722
$ws.Range("T180").Value = 18

# New row 181: Blue Giant / Segunda
$ws.Range("A181").Value = 10
$ws.Range("B181").Value = "Vega Modelo de Temuco"
$ws.Range("C181").Value = "La Araucanía"
$ws.Range("D181").Value = 44610
$ws.Range("E181").Value = 9
$ws.Range("F181").Value = "Fruta"
$ws.Range("G181").Value = 100103
$ws.Range("H181").Value = "Frutos de hueso (carozo)"
$ws.Range("I181").Value = 100103002
$ws.Range("J181").Value = "Ciruela"
$ws.Range("K181").Value = "Blue Giant"
$ws.Range("L181").Value = "Segunda"
$ws.Range("M181").Value = 2
$ws.Range("N181").Value = 160000
$ws.Range("O181").Value = 160000
$ws.Range("P181").Value = 160000
$ws.Range("Q181").Value = "$/bins (450 kilos)"
$ws.Range("R181").Value = "Región de O'Higgins"
$ws.Range("S181").Value = 356
$ws.Range("T181").Value = 450
